$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 and 17: Avalanche/WrappedBTC swap with new data
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.73"
$ws.Range("E16").Value = "  -6.11%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "64.372.99"
$ws.Range("E17").Value = "  +0.56%  "

# Row 2
$ws.Range("D2").Value = "64.300.72"
$ws.Range("E2").Value = "  +0.69%  "

# Row 3
$ws.Range("D3").Value = "3.496.75"
$ws.Range("E3").Value = "  -0.24%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.79"
$ws.Range("E5").Value = "  +0.60%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.19"
$ws.Range("E6").Value = "  +1.07%  "

# Row 7
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("E8").Value = "  +0.61%  "

# Row 9
$ws.Range("E9").Value = "  +0.51%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.29"
$ws.Range("E10").Value = "  +2.65%  "

# Row 11
$ws.Range("E11").Value = "  +2.82%  "

# Row 12
$ws.Range("D12").Value = "4.092.76"
$ws.Range("E12").Value = "  -0.50%  "

# Row 13
$ws.Range("E13").Value = "  +1.24%  "

# Row 14
$ws.Range("E14").Value = "  +1.61%  "

# Row 15
$ws.Range("D15").Value = "3.495.72"
$ws.Range("E15").Value = "  -0.41%  "

# Row 18
$ws.Range("E18").Value = "  +0.70%  "

# Row 19
$ws.Range("E19").Value = "  +2.82%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.54"
$ws.Range("E20").Value = "  -2.53%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "393.34"
$ws.Range("E21").Value = "  +2.86%  "

# Row 22
$ws.Range("E22").Value = "  +0.14%  "

# Row 23
$ws.Range("D23").Value = "3.636.25"
$ws.Range("E23").Value = "  -0.46%  "

# Row 24
$ws.Range("E24").Value = "  +1.22%  "

# Row 25
$ws.Range("E25").Value = "  +0.02%  "

# Row 26
$ws.Range("E26").Value = "  +1.52%  "

# Row 27
$ws.Range("E27").Value = "  +1.55%  "

# Row 28
$ws.Range("E28").Value = "  -0.20%  "

# Row 29
$ws.Range("E29").Value = "  -1.17%  "

# Row 30
$ws.Range("E30").Value = "  +0.84%  "

# Row 31
$ws.Range("E31").Value = "  -1.58%  "

# Row 32
$ws.Range("E32").Value = "  -5.63%  "

# Row 33
$ws.Range("D33").Value = "3.519.28"
$ws.Range("E33").Value = "  -0.03%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.152"
$ws.Range("E34").Value = "  +4.89%  "

# Row 35
$ws.Range("E35").Value = "  +0.05%  "

# Row 36
$ws.Range("E36").Value = "  -0.12%  "

# Row 37
$ws.Range("E37").Value = "  -3.57%  "

# Row 38
$ws.Range("E38").Value = "  -0.09%  "

# Row 39
$ws.Range("E39").Value = "  -0.69%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "167.46"
$ws.Range("E40").Value = "  +4.18%  "

# Row 41
$ws.Range("E41").Value = "  -0.48%  "

# Row 42
$ws.Range("E42").Value = "  -0.26%  "

# Row 43
$ws.Range("E43").Value = "  -0.19%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.44"
$ws.Range("E44").Value = "  -3.79%  "

# Row 45
$ws.Range("E45").Value = "  +0.05%  "

# Row 46
$ws.Range("E46").Value = "  +2.71%  "

# Row 47
$ws.Range("E47").Value = "  -3.74%  "

# Row 48
$ws.Range("D48").Value = "2.468.84"
$ws.Range("E48").Value = "  -0.47%  "

# Row 49
$ws.Range("E49").Value = "  -0.44%  "

# Row 50
$ws.Range("E50").Value = "  -1.33%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0259"
$ws.Range("E51").Value = "  -0.87%  "
